$d = $word.ActiveDocument

$replacements = @(
    @("33÷5=6, 3", "34÷8=4, 2"),
    @("63÷4=15, 3", "74÷8=9, 2"),
    @("83÷8=10, 3", "75÷9=8, 3"),
    @("26÷2=13, 0", "13÷7=1, 6"),
    @("86÷3=28, 2", "58÷6=9, 4"),
    @("97÷2=48, 1", "68÷9=7, 5"),
    @("37÷9=4, 1", "50÷8=6, 2"),
    @("74÷9=8, 2", "61÷8=7, 5"),
    @("51÷6=8, 3", "99÷5=19, 4"),
    @("55÷6=9, 1", "91÷3=30, 1"),
    @("42÷9=4, 6", "40÷4=10, 0"),
    @("33÷3=11, 0", "96÷4=24, 0"),
    @("32÷2=16, 0", "27÷9=3, 0"),
    @("22÷8=2, 6", "46÷2=23, 0"),
    @("25÷7=3, 4", "13÷7=1, 6"),
    @("80÷8=10, 0", "30÷9=3, 3"),
    @("59÷5=11, 4", "59÷6=9, 5"),
    @("39÷3=13, 0", "36÷2=18, 0"),
    @("52÷9=5, 7", "89÷8=11, 1"),
    @("72÷3=24, 0", "67÷4=16, 3"),
    @("12÷3=4, 0", "35÷9=3, 8"),
    @("37÷8=4, 5", "52÷2=26, 0"),
    @("85÷2=42, 1", "63÷4=15, 3"),
    @("77÷2=38, 1", "95÷3=31, 2"),
    @("61÷6=10, 1", "72÷2=36, 0")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
